# Edit slide 4 ("Install MongoDB (server) and Robo3T (client)") in the
# Text Placeholder shape: update the first video title/description and
# the YouTube URL beneath it.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1: replace the single run with three runs -----------------
# "How to install MongoDB 4.4.1 on Windows 10 (2020)"
#   -> "How To Install MongoDB On Windows 10 | MongoDB Installation | MongoDB Tutorial | "
#      + "Simplilearn" + " (2021)"
$para1 = $tr.Paragraphs(1)
$run1  = $para1.Runs(1)

# Re-set the run's own text (instead of the paragraph's) so the engine
# performs a plain replacement instead of diffing against the old text
# and splitting off a shared prefix into its own run.
$run1.Text = "How To Install MongoDB On Windows 10 | MongoDB Installation | MongoDB Tutorial | "

# InsertAfter clones the formatting (rPr) of the preceding run, giving us
# the "Avenir Medium" / sz=2000 formatting on the new runs automatically.
$run2 = $run1.InsertAfter("Simplilearn")
$run3 = $run2.InsertAfter(" (2021)")

# --- Paragraph 2: swap the YouTube video id in the URL text --------------
# The hyperlink relationship (rId2) itself is unchanged; only the visible
# text of the run changes.
$para2 = $tr.Paragraphs(2)
$urlRun = $para2.Runs(1)
$urlRun.Text = "https://www.youtube.com/watch?v=Z478ODY4ceQ"
